$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Automatic tracker update: settle pending results for concluded matches,
# fix the shifted/duplicated rows around 54-63, and append newly scraped
# fixtures (rows 64-72).

# --- Rows where only columns A-F change (results not yet settled; G/H stay blank) ---
$ws.Cells.Item(56, 1).Value = 14349723
$ws.Cells.Item(56, 2).NumberFormat = '@'
$ws.Cells.Item(56, 2).Value = '2025-08-05'
$ws.Cells.Item(56, 2).Style = 'Normal'
$ws.Cells.Item(56, 3).Value = 'Gonzalo Bueno'
$ws.Cells.Item(56, 4).Value = 'Ryan Nijboer'
$ws.Cells.Item(56, 5).Value = 'Gana Ryan Nijboer'
$ws.Cells.Item(56, 6).Value = 2.75

$ws.Cells.Item(59, 1).Value = 14349616
$ws.Cells.Item(59, 2).NumberFormat = '@'
$ws.Cells.Item(59, 2).Value = '2025-08-05'
$ws.Cells.Item(59, 2).Style = 'Normal'
$ws.Cells.Item(59, 3).Value = 'Abdullah Shelbayh'
$ws.Cells.Item(59, 4).Value = 'Ugo Blanchet'
$ws.Cells.Item(59, 5).Value = 'Gana Abdullah Shelbayh'
$ws.Cells.Item(59, 6).Value = 2.75

$ws.Cells.Item(60, 1).Value = 14349615
$ws.Cells.Item(60, 2).NumberFormat = '@'
$ws.Cells.Item(60, 2).Value = '2025-08-05'
$ws.Cells.Item(60, 2).Style = 'Normal'
$ws.Cells.Item(60, 3).Value = 'Rafael Jodar'
$ws.Cells.Item(60, 4).Value = 'Marc-Andrea Huesler'
$ws.Cells.Item(60, 5).Value = 'Gana Marc-Andrea Huesler'
$ws.Cells.Item(60, 6).Value = 1.91

$ws.Cells.Item(61, 1).Value = 14349721
$ws.Cells.Item(61, 2).NumberFormat = '@'
$ws.Cells.Item(61, 2).Value = '2025-08-05'
$ws.Cells.Item(61, 2).Style = 'Normal'
$ws.Cells.Item(61, 3).Value = 'Max Alcala Gurri'
$ws.Cells.Item(61, 4).Value = 'Lorenzo Giustino'
$ws.Cells.Item(61, 5).Value = 'Gana Lorenzo Giustino'
$ws.Cells.Item(61, 6).Value = 2.1

$ws.Cells.Item(62, 1).Value = 14349731
$ws.Cells.Item(62, 2).NumberFormat = '@'
$ws.Cells.Item(62, 2).Value = '2025-08-05'
$ws.Cells.Item(62, 2).Style = 'Normal'
$ws.Cells.Item(62, 3).Value = 'Lukas Neumayer'
$ws.Cells.Item(62, 4).Value = 'Mariano Kestelboim'
$ws.Cells.Item(62, 5).Value = 'Gana Mariano Kestelboim'
$ws.Cells.Item(62, 6).Value = 3.75

$ws.Cells.Item(63, 1).Value = 14350800
$ws.Cells.Item(63, 2).NumberFormat = '@'
$ws.Cells.Item(63, 2).Value = '2025-08-05'
$ws.Cells.Item(63, 2).Style = 'Normal'
$ws.Cells.Item(63, 3).Value = 'Chun-Hsin Tseng'
$ws.Cells.Item(63, 4).Value = 'Zachary Svajda'
$ws.Cells.Item(63, 5).Value = 'Gana Chun-Hsin Tseng'
$ws.Cells.Item(63, 6).Value = 4

# --- Brand new rows appended at the bottom (results not yet settled) ---
$ws.Cells.Item(64, 1).Value = 14350770
$ws.Cells.Item(64, 2).NumberFormat = '@'
$ws.Cells.Item(64, 2).Value = '2025-08-05'
$ws.Cells.Item(64, 2).Style = 'Normal'
$ws.Cells.Item(64, 3).Value = 'Valentin Royer'
$ws.Cells.Item(64, 4).Value = 'Juan Pablo Ficovich'
$ws.Cells.Item(64, 5).Value = 'Gana Juan Pablo Ficovich'
$ws.Cells.Item(64, 6).Value = 2.62

$ws.Cells.Item(65, 1).Value = 14350777
$ws.Cells.Item(65, 2).NumberFormat = '@'
$ws.Cells.Item(65, 2).Value = '2025-08-05'
$ws.Cells.Item(65, 2).Style = 'Normal'
$ws.Cells.Item(65, 3).Value = 'Juan Manuel Cerundolo'
$ws.Cells.Item(65, 4).Value = 'Coleman Wong'
$ws.Cells.Item(65, 5).Value = 'Gana Juan Manuel Cerundolo'
$ws.Cells.Item(65, 6).Value = 2.25

$ws.Cells.Item(66, 1).Value = 14350783
$ws.Cells.Item(66, 2).NumberFormat = '@'
$ws.Cells.Item(66, 2).Value = '2025-08-05'
$ws.Cells.Item(66, 2).Style = 'Normal'
$ws.Cells.Item(66, 3).Value = 'Christopher Eubanks'
$ws.Cells.Item(66, 4).Value = 'Nikoloz Basilashvili'
$ws.Cells.Item(66, 5).Value = 'Gana Christopher Eubanks'
$ws.Cells.Item(66, 6).Value = 2.25

$ws.Cells.Item(67, 1).Value = 14350794
$ws.Cells.Item(67, 2).NumberFormat = '@'
$ws.Cells.Item(67, 2).Value = '2025-08-05'
$ws.Cells.Item(67, 2).Style = 'Normal'
$ws.Cells.Item(67, 3).Value = 'Omar Jasika'
$ws.Cells.Item(67, 4).Value = 'Terence Atmane'
$ws.Cells.Item(67, 5).Value = 'Gana Omar Jasika'
$ws.Cells.Item(67, 6).Value = 3.5

$ws.Cells.Item(68, 1).Value = 14350771
$ws.Cells.Item(68, 2).NumberFormat = '@'
$ws.Cells.Item(68, 2).Value = '2025-08-05'
$ws.Cells.Item(68, 2).Style = 'Normal'
$ws.Cells.Item(68, 3).Value = 'Adrian Mannarino'
$ws.Cells.Item(68, 4).Value = 'Mitchell Krueger'
$ws.Cells.Item(68, 5).Value = 'Gana Mitchell Krueger'
$ws.Cells.Item(68, 6).Value = 2.75

$ws.Cells.Item(69, 1).Value = 14350776
$ws.Cells.Item(69, 2).NumberFormat = '@'
$ws.Cells.Item(69, 2).Value = '2025-08-05'
$ws.Cells.Item(69, 2).Style = 'Normal'
$ws.Cells.Item(69, 3).Value = 'Aleksandar Vukic'
$ws.Cells.Item(69, 4).Value = 'Murphy Cassone'
$ws.Cells.Item(69, 5).Value = 'Gana Murphy Cassone'
$ws.Cells.Item(69, 6).Value = 3.4

$ws.Cells.Item(70, 1).Value = 14350772
$ws.Cells.Item(70, 2).NumberFormat = '@'
$ws.Cells.Item(70, 2).Value = '2025-08-05'
$ws.Cells.Item(70, 2).Style = 'Normal'
$ws.Cells.Item(70, 3).Value = 'Patrick Kypson'
$ws.Cells.Item(70, 4).Value = 'Matteo Gigante'
$ws.Cells.Item(70, 5).Value = 'Gana Patrick Kypson'
$ws.Cells.Item(70, 6).Value = 2.25

$ws.Cells.Item(71, 1).Value = 14350934
$ws.Cells.Item(71, 2).NumberFormat = '@'
$ws.Cells.Item(71, 2).Value = '2025-08-05'
$ws.Cells.Item(71, 2).Style = 'Normal'
$ws.Cells.Item(71, 3).Value = 'Hanyu Guo'
$ws.Cells.Item(71, 4).Value = 'Anca Todoni'
$ws.Cells.Item(71, 5).Value = 'Gana Hanyu Guo'
$ws.Cells.Item(71, 6).Value = 2.5

$ws.Cells.Item(72, 1).Value = 14349601
$ws.Cells.Item(72, 2).NumberFormat = '@'
$ws.Cells.Item(72, 2).Value = '2025-08-05'
$ws.Cells.Item(72, 2).Style = 'Normal'
$ws.Cells.Item(72, 3).Value = 'Martin Krumich'
$ws.Cells.Item(72, 4).Value = 'Patrick Zahraj'
$ws.Cells.Item(72, 5).Value = 'Gana Patrick Zahraj'
$ws.Cells.Item(72, 6).Value = 2.38

# --- Rows where A-F change AND results (G/H) are now known ---
$ws.Cells.Item(29, 1).Value = 14339487
$ws.Cells.Item(29, 2).NumberFormat = '@'
$ws.Cells.Item(29, 2).Value = '2025-08-04'
$ws.Cells.Item(29, 2).Style = 'Normal'
$ws.Cells.Item(29, 3).Value = 'Santiago Rodriguez Taverna'
$ws.Cells.Item(29, 4).Value = 'Nikolas Sanchez Izquierdo'
$ws.Cells.Item(29, 5).Value = 'Gana Santiago Rodriguez Taverna'
$ws.Cells.Item(29, 6).Value = 1.83
$ws.Cells.Item(29, 7).Value = 'Acierto'
$ws.Cells.Item(29, 8).Value = 0.83

$ws.Cells.Item(31, 1).Value = 14339502
$ws.Cells.Item(31, 2).NumberFormat = '@'
$ws.Cells.Item(31, 2).Value = '2025-08-04'
$ws.Cells.Item(31, 2).Style = 'Normal'
$ws.Cells.Item(31, 3).Value = 'Jan Choinski'
$ws.Cells.Item(31, 4).Value = 'Geoffrey Blancaneaux'
$ws.Cells.Item(31, 5).Value = 'Gana Geoffrey Blancaneaux'
$ws.Cells.Item(31, 6).Value = 2.63
$ws.Cells.Item(31, 7).Value = 'Acierto'
$ws.Cells.Item(31, 8).Value = 1.63

$ws.Cells.Item(51, 1).Value = 14349617
$ws.Cells.Item(51, 2).NumberFormat = '@'
$ws.Cells.Item(51, 2).Value = '2025-08-05'
$ws.Cells.Item(51, 2).Style = 'Normal'
$ws.Cells.Item(51, 3).Value = 'Daniil Glinka'
$ws.Cells.Item(51, 4).Value = 'Radu Albot'
$ws.Cells.Item(51, 5).Value = 'Gana Radu Albot'
$ws.Cells.Item(51, 6).Value = 1.83
$ws.Cells.Item(51, 7).Value = 'Fallo'
$ws.Cells.Item(51, 8).Value = -1

$ws.Cells.Item(52, 1).Value = 14349604
$ws.Cells.Item(52, 2).NumberFormat = '@'
$ws.Cells.Item(52, 2).Value = '2025-08-05'
$ws.Cells.Item(52, 2).Style = 'Normal'
$ws.Cells.Item(52, 3).Value = 'Gonzalo Villanueva'
$ws.Cells.Item(52, 4).Value = 'Elmer Moller'
$ws.Cells.Item(52, 5).Value = 'Gana Gonzalo Villanueva'
$ws.Cells.Item(52, 6).Value = 5.5
$ws.Cells.Item(52, 7).Value = 'Fallo'
$ws.Cells.Item(52, 8).Value = -1

$ws.Cells.Item(54, 1).Value = 14349613
$ws.Cells.Item(54, 2).NumberFormat = '@'
$ws.Cells.Item(54, 2).Value = '2025-08-05'
$ws.Cells.Item(54, 2).Style = 'Normal'
$ws.Cells.Item(54, 3).Value = 'Francesco Maestrelli'
$ws.Cells.Item(54, 4).Value = 'Maximus Jones'
$ws.Cells.Item(54, 5).Value = 'Gana Maximus Jones'
$ws.Cells.Item(54, 6).Value = 3.5
$ws.Cells.Item(54, 7).Value = 'Fallo'
$ws.Cells.Item(54, 8).Value = -1

$ws.Cells.Item(55, 1).Value = 14349730
$ws.Cells.Item(55, 2).NumberFormat = '@'
$ws.Cells.Item(55, 2).Value = '2025-08-05'
$ws.Cells.Item(55, 2).Style = 'Normal'
$ws.Cells.Item(55, 3).Value = 'Zdenek Kolar'
$ws.Cells.Item(55, 4).Value = 'Zsombor Piros'
$ws.Cells.Item(55, 5).Value = 'Gana Zdenek Kolar'
$ws.Cells.Item(55, 6).Value = 3.4
$ws.Cells.Item(55, 7).Value = 'Fallo'
$ws.Cells.Item(55, 8).Value = -1

$ws.Cells.Item(57, 1).Value = 14349602
$ws.Cells.Item(57, 2).NumberFormat = '@'
$ws.Cells.Item(57, 2).Value = '2025-08-05'
$ws.Cells.Item(57, 2).Style = 'Normal'
$ws.Cells.Item(57, 3).Value = 'Mats Rosenkranz'
$ws.Cells.Item(57, 4).Value = 'Clement Tabur'
$ws.Cells.Item(57, 5).Value = 'Gana Mats Rosenkranz'
$ws.Cells.Item(57, 6).Value = 3.25
$ws.Cells.Item(57, 7).Value = 'Acierto'
$ws.Cells.Item(57, 8).Value = 2.25

$ws.Cells.Item(58, 1).Value = 14339497
$ws.Cells.Item(58, 2).NumberFormat = '@'
$ws.Cells.Item(58, 2).Value = '2025-08-05'
$ws.Cells.Item(58, 2).Style = 'Normal'
$ws.Cells.Item(58, 3).Value = 'Justin Engel'
$ws.Cells.Item(58, 4).Value = 'Alejandro Moro Canas'
$ws.Cells.Item(58, 5).Value = 'Gana Justin Engel'
$ws.Cells.Item(58, 6).Value = 3
$ws.Cells.Item(58, 7).Value = 'Acierto'
$ws.Cells.Item(58, 8).Value = 2

